$d = $word.ActiveDocument

# Locate the paragraph that credits Jan Hollan / links to the astromap archive.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Jan Hollan*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    # Exclude the trailing paragraph mark from the range we rewrite.
    $r.End = $r.End - 1
    $r.Delete()
    $r.InsertAfter("de Jan Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).")
}
